$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = "Medium (Requires Calling Gradient Calculation Object For Every Activation Function Calculation Object.)"
$ws.Range("E16").Value = "Hard (Nested Arrays And Functions Inside Function Parameters.)"
$ws.Range("F14").Select()
